# Insert a new "snapshot" column before the trailing "nom"/"url_produit"
# columns (previously FL, FM), shifting them right to FM/FN.
#
# Column FK (167) held the most recent price snapshot; the new column FL
# (168) becomes the newest snapshot: its header (row 1) gets the new
# timestamp, and each product row gets the same price that was already in
# FK (carried forward), or stays blank if FK itself was blank for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLastCol = 167   # FK - last existing price-history column before nom/url_produit
$newCol     = 168   # FL - newly inserted price-history column

# Shift "nom" / "url_produit" (and everything else) one column to the right,
# inserting a fresh blank column at FL.
$ws.Columns.Item($newCol).Insert()

# New header timestamp for the inserted column.
$ws.Cells.Item(1, $newCol).Value = "2026-02-04 17:37:03"

# Carry the last known price forward into the new column for every data row.
$lastRow = 208
for ($r = 2; $r -le $lastRow; $r++) {
    $prevVal = $ws.Cells.Item($r, $oldLastCol).Value()
    if ($null -ne $prevVal -and $prevVal -ne "") {
        $ws.Cells.Item($r, $newCol).Value = $prevVal
    }
}
